$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column cells that would otherwise look like a plain decimal number
# need their number format forced to Text first, matching the original
# inline-string (text) storage so Excel does not coerce "0.1700" -> 0.17 etc.
$textCells = @("D4","D5","D6","D8","D9","D10","D11","D13","D14","D15","D16","D17","D20","D22","D23","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.810.25"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.856.21"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  -1.61%  "
$ws.Range("D5").Value = "319.49"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").Value = "0.3755"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "0.07354"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "0.8777"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "21.59"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.866.02"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "6.764"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "5.453"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "0.07144"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "89.11"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "27.829.88"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "5.222"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").Value = "2.080.42"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "1.984"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "155.42"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "2.173"
$ws.Range("E28").Value = "  +9.45%  "
$ws.Range("D29").Value = "5.373"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "119.05"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "0.7778"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "4.552"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "2.929"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "1.133"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01986"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05352"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.904"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.166"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("D42").Value = "0.1700"
$ws.Range("D43").Value = "0.5150"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "8.854"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "10.75"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "108.31"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "0.4773"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "0.06479"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "1.693"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "1.012"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "1.852"
$ws.Range("E51").Value = "  -2.35%  "
